# Magalie virospore quant - add raw_data_joy sheet with updated figures,
# tweak raw_data sheet view/column width, matching the commit:
# "Updated figures and included new code for quantifying virospores for Magalie revisions"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Tweak the "raw_data" sheet: widen column F, and change the
#    selection to the full data range without it being the active tab.
# ---------------------------------------------------------------------
$rawData = $wb.Worksheets.Item("raw_data")
$rawData.Columns.Item(6).ColumnWidth = 21
[void]$rawData.Range("A1:N7").Select()

# ---------------------------------------------------------------------
# 2) Insert a new worksheet "raw_data_joy" right after "raw_data" (and
#    before "raw_data_calculation"), carrying Joy's updated figures.
# ---------------------------------------------------------------------
$rawDataJoy = $wb.Worksheets.Add($null, $rawData)
$rawDataJoy.Name = "raw_data_joy"

$headers = @("Plaque", "PFU_Diameter_cm", "Sample", "Region", "VPFU", "Free_PFU_10e3", "Spores_CFU_10e3", "Region_start_fromPFUcenter_cm", "Region_end-fromPFUcenter_cm", "Distance_from_PFU_center_cm")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $rawDataJoy.Cells.Item(1, $col + 1).Value = $headers[$col]
}

$rows = @(
    @("A", 0.2, 1.1, "Center",  23, 54, 0, 0,    0.05, 0),
    @("A", 0.2, 1.2, "Center",  42, 18, 0, 0,    0.05, 0.025),
    @("A", 0.2, 2.1, "Annulus", 292, 15, 0, 0.05, 0.1, 0.05),
    @("A", 0.2, 2.2, "Annulus", 602, 28, 0, 0.05, 0.1, 0.075),
    @("A", 0.2, 3.1, "Lawn",    29, 2, 2, 0.1,  0.2, 0.1),
    @("A", 0.2, 3.2, "Lawn",    77, 2, 3, 0.1,  0.2, 0.2)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowValues = $rows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $rawDataJoy.Cells.Item($r + 2, $c + 1).Value = $rowValues[$c]
    }
}

# Column widths for the new sheet
$rawDataJoy.Columns.Item(6).ColumnWidth = 20
$rawDataJoy.Columns.Item(7).ColumnWidth = 21.5
$rawDataJoy.Columns.Item(10).ColumnWidth = 22.5

# Select whole column H (to match the author's "tabSelected" snapshot) and
# make this new sheet the active tab of the workbook.
[void]$rawDataJoy.Range("H1:H1048576").Select()
[void]$rawDataJoy.Activate()

Write-Host "raw_data_joy sheet created with updated figures"
